$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Purchase 22-23")
$ws2 = $wb.Worksheets.Item("Sale 22-23")

# --- Formatting fix-ups on rows that keep their row number -------------------
# Row 24 takes on the "section boundary" look that row 17 already has
# (thin border without a bottom rule - style group 26/27/28 in the original file).
$ws1.Range("A17:F17").Copy()
$ws1.Range("A24:F24").PasteSpecial(-4122)

# Row 25's F cell reverts to the plain (non-boundary) look used elsewhere in the block.
$ws1.Range("F23").Copy()
$ws1.Range("F25").PasteSpecial(-4122)

# --- Insert the new transaction row ------------------------------------------
# A new row is inserted at 26 (pushing the "2nd party" block down by one row,
# preserving the existing blank spacer rows).
$ws1.Rows("26:26").Insert()

# Give the new row the same "section boundary" look the old row 27
# (now shifted to row 29) already carries.
$ws1.Range("A29:F29").Copy()
$ws1.Range("A26:F26").PasteSpecial(-4122)

# --- New row 26 data -----------------------------------------------------------
$ws1.Range("B26").Value = 45180
$ws1.Range("C26").Value = "50/23-24"
$ws1.Range("D26").Value = "Namrata Rubber Product Pvt Ltd"
$ws1.Range("E26").Value = 33040
$ws1.Range("F26").Formula = "=F25+E26"

$excel.CutCopyMode = 0

# --- View / selection state ---------------------------------------------------
$ws2.Range("D20").Select()
$ws1.Activate()
$ws1.Range("G20").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1

$wb.Application.Calculate()
